$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows 386-464 (update through 2021-12-08 / "aggiornamento fino a 8/12")
# Columns: A = date serial, B = nuovi pos., C = somma mobile 7gg., D = somma mobile 7gg. per 100mila abitanti
$data = @(
    ,@(44460, 0, 2, 12.52348152786475)
    ,@(44461, 0, 2, 12.52348152786475)
    ,@(44462, 2, 4, 25.04696305572949)
    ,@(44463, 2, 6, 37.57044458359425)
    ,@(44464, 0, 5, 31.30870381966186)
    ,@(44465, 0, 4, 25.04696305572949)
    ,@(44466, 0, 4, 25.04696305572949)
    ,@(44467, 0, 4, 25.04696305572949)
    ,@(44468, 0, 4, 25.04696305572949)
    ,@(44469, 0, 2, 12.52348152786475)
    ,@(44470, 1, 1, 6.261740763932373)
    ,@(44471, 0, 1, 6.261740763932373)
    ,@(44472, 0, 1, 6.261740763932373)
    ,@(44473, 0, 1, 6.261740763932373)
    ,@(44474, 0, 1, 6.261740763932373)
    ,@(44475, 1, 2, 12.52348152786475)
    ,@(44476, 5, 7, 43.83218534752661)
    ,@(44477, 2, 8, 50.09392611145898)
    ,@(44478, 0, 8, 50.09392611145898)
    ,@(44479, 0, 8, 50.09392611145898)
    ,@(44480, 1, 9, 56.35566687539136)
    ,@(44481, 0, 9, 56.35566687539136)
    ,@(44482, 0, 8, 50.09392611145898)
    ,@(44483, 1, 4, 25.04696305572949)
    ,@(44484, 0, 2, 12.52348152786475)
    ,@(44485, 0, 2, 12.52348152786475)
    ,@(44486, 1, 3, 18.78522229179712)
    ,@(44487, 0, 2, 12.52348152786475)
    ,@(44488, 0, 2, 12.52348152786475)
    ,@(44489, 0, 2, 12.52348152786475)
    ,@(44490, 0, 1, 6.261740763932373)
    ,@(44491, 0, 1, 6.261740763932373)
    ,@(44492, 0, 1, 6.261740763932373)
    ,@(44493, 1, 1, 6.261740763932373)
    ,@(44494, 0, 1, 6.261740763932373)
    ,@(44495, 0, 1, 6.261740763932373)
    ,@(44496, 3, 4, 25.04696305572949)
    ,@(44497, 0, 4, 25.04696305572949)
    ,@(44498, 0, 4, 25.04696305572949)
    ,@(44499, 0, 4, 25.04696305572949)
    ,@(44500, 0, 3, 18.78522229179712)
    ,@(44501, 0, 3, 18.78522229179712)
    ,@(44502, 0, 3, 18.78522229179712)
    ,@(44503, 0, 0, 0)
    ,@(44504, 0, 0, 0)
    ,@(44505, 1, 1, 6.261740763932373)
    ,@(44506, 0, 1, 6.261740763932373)
    ,@(44507, 0, 1, 6.261740763932373)
    ,@(44508, 0, 1, 6.261740763932373)
    ,@(44509, 0, 1, 6.261740763932373)
    ,@(44510, 0, 1, 6.261740763932373)
    ,@(44511, 1, 2, 12.52348152786475)
    ,@(44512, 0, 1, 6.261740763932373)
    ,@(44513, 0, 1, 6.261740763932373)
    ,@(44514, 0, 1, 6.261740763932373)
    ,@(44515, 1, 2, 12.52348152786475)
    ,@(44516, 9, 11, 68.8791484032561)
    ,@(44517, 0, 11, 68.8791484032561)
    ,@(44518, 0, 10, 62.61740763932373)
    ,@(44519, 4, 14, 87.66437069505322)
    ,@(44520, 1, 15, 93.9261114589856)
    ,@(44521, 1, 16, 100.187852222918)
    ,@(44522, 7, 22, 137.7582968065122)
    ,@(44523, 0, 13, 81.40262993112086)
    ,@(44524, 21, 34, 212.8991859737007)
    ,@(44525, 2, 36, 225.4226675015654)
    ,@(44526, 3, 35, 219.160926737633)
    ,@(44527, 4, 38, 237.9461490294302)
    ,@(44528, 10, 47, 294.3018159048215)
    ,@(44529, 9, 49, 306.8252974326863)
    ,@(44530, 1, 50, 313.0870381966187)
    ,@(44531, 1, 30, 187.8522229179712)
    ,@(44532, 2, 30, 187.8522229179712)
    ,@(44533, 12, 39, 244.2078897933626)
    ,@(44534, 1, 36, 225.4226675015654)
    ,@(44535, 0, 26, 162.8052598622417)
    ,@(44536, 10, 27, 169.0670006261741)
    ,@(44537, 0, 26, 162.8052598622417)
    ,@(44538, 3, 28, 175.3287413901064)
)

$startRow = 386
$endRow = $startRow + $data.Count - 1

# Extend column A date formatting (style of A385) down through the new rows in one shot
$ws.Range("A385").Copy() | Out-Null
$ws.Range("A{0}:A{1}" -f $startRow, $endRow).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}